$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("id") holds purely-numeric-looking strings in the source file
# (stored as text/inlineStr there). Force text number format first so
# the COM layer keeps them as strings instead of coercing to numbers.
# (Row 4 "phone" values contain dashes, so they naturally stay text.)
$ws.Range("B2:G2").NumberFormat = "@"

$ws.Range("B2").Value = "836942"
$ws.Range("C2").Value = "727809"
$ws.Range("D2").Value = "949621"
$ws.Range("E2").Value = "573574"
$ws.Range("F2").Value = "733386"
$ws.Range("G2").Value = "803444"

$ws.Range("B4").Value = "522-0-1"
$ws.Range("C4").Value = "837-85-44"
$ws.Range("D4").Value = "407-81-47"
$ws.Range("E4").Value = "324-24-64"
$ws.Range("F4").Value = "400-96-1"
$ws.Range("G4").Value = "661-43-96"
